# Fruta / hortaliza, semanal
#
# Insert six new weekly price records for Kiwi (Femacal de La Calera,
# "Región de O'Higgins") right after the existing row 627, pushing every
# subsequent row down by six. The sheet's used range grows from A1:T727
# to A1:T733.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the six new records (rows 628-633); everything that used
# to live at 628..727 slides down to 634..733.
$ws.Rows("628:633").Insert()

$newRows = @(
  @(3,"Femacal de La Calera","Coquimbo",44776,5,"Fruta",100101,"Berries",100101007,"Kiwi","Hayward","Especial",75,7000,7000,7000,"`$/bandeja 10 kilos","Región de O'Higgins",700,10),
  @(3,"Femacal de La Calera","Coquimbo",44776,5,"Fruta",100101,"Berries",100101007,"Kiwi","Hayward","Especial",56,10000,10000,10000,"`$/caja 15 kilos","Región de O'Higgins",667,15),
  @(3,"Femacal de La Calera","Coquimbo",44776,5,"Fruta",100101,"Berries",100101007,"Kiwi","Hayward","Primera",70,6000,6000,6000,"`$/bandeja 10 kilos","Región de O'Higgins",600,10),
  @(3,"Femacal de La Calera","Coquimbo",44776,5,"Fruta",100101,"Berries",100101007,"Kiwi","Hayward","Primera",67,9000,9000,9000,"`$/caja 15 kilos","Región de O'Higgins",600,15),
  @(3,"Femacal de La Calera","Coquimbo",44776,5,"Fruta",100101,"Berries",100101007,"Kiwi","Hayward","Segunda",80,5000,5000,5000,"`$/bandeja 10 kilos","Región de O'Higgins",500,10),
  @(3,"Femacal de La Calera","Coquimbo",44776,5,"Fruta",100101,"Berries",100101007,"Kiwi","Hayward","Segunda",60,7500,7500,7500,"`$/caja 15 kilos","Región de O'Higgins",500,15)
)

$startRow = 628
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowData = $newRows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $rowData[$c - 1]
    }
}
